# Small adjustments to consumption. Added elysis SE03. Adjusted offshore wind SE.
#
# - Un-hide the previously filtered "2040" rows (7-14, 16-20) by clearing the
#   autofilter criteria (was filtered to Scenario = "National Trends").
# - Append two new demand rows (SE03 / ES00, hydrogen, Distributed Energy, 2040)
# - Re-apply the autofilter over the full data range with no filter criteria.
# - Move the active selection down to the new last cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the existing "National Trends" column filter and unhide rows ---
$ws.AutoFilterMode = $false
foreach ($r in 7..20) {
    $ws.Rows.Item($r).Hidden = $false
}

# Re-apply the autofilter across the original range (A1:E26), no filter
# criteria. Do this BEFORE appending the new rows below so the filter range
# doesn't auto-grow to include them.
[void]$ws.Range("A1:E26").AutoFilter()

# Keep the workbook's hidden _FilterDatabase defined name in sync with the
# new autofilter range (Excel normally maintains this automatically).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "demand!_FilterDatabase") {
        $n.RefersTo = "=demand!`$A`$1:`$E`$26"
    }
}

# --- Append the two new rows (this also grows the sheet dimension to A1:E28) ---
$ws.Range("A27").Value = "SE03"
$ws.Range("B27").Value = "hydrogen"
$ws.Range("C27").Value = "Distributed Energy"
$ws.Range("D27").Value = 2040
$ws.Range("E27").Value = 600

$ws.Range("A28").Value = "ES00"
$ws.Range("B28").Value = "hydrogen"
$ws.Range("C28").Value = "Distributed Energy"
$ws.Range("D28").Value = 2040
$ws.Range("E28").Value = 1000

# Match the formatting used by the rest of the data rows.
$ws.Range("A27:E28").NumberFormat = $ws.Range("A26:E26").NumberFormat

# Move the selection to follow the newly-added last row.
[void]$ws.Range("E29").Select()
